# edit.ps1 - applies the "added SSL depricating notice" commit to the deck.
#
# Summary of changes (see diff):
#   1. Fix the PGP bullet's "assymmetric" typo -> "asymmetric" and collapse
#      the three runs that spelled it out into the single merged run the
#      author ended up with.
#   2. Bold + color (red) the "a successor of SSL 3.0" phrase inside the
#      TLS bullet, splitting that run into three runs.
#   3. Add a new standalone red/bold textbox reading
#      "SSL was deprecated in 2015" near the SSL/TLS diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The bullet list lives in shape #2 on the slide ("TextBox 2" / id 3).
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- 1. PGP bullet: fix typo, merge the runs into one --------------------
$full = $tr.Text
$oldPgp = " = Pretty Good Privacy, since 1991, uses both symmetric an assymmetric encryption, the de facto standard for email security."
$newPgp = " = Pretty Good Privacy, since 1991, uses both symmetric an asymmetric encryption, the de facto standard for email security."
$idxPgp = $full.IndexOf($oldPgp)
$rngPgp = $tr.Characters($idxPgp + 1, $oldPgp.Length)
$rngPgp.Text = $newPgp

# Editing the run text makes this auto-fit ("shrink/grow shape to fit
# text") shape recompute its height against live layout metrics; put the
# shape back to its authored height (6555641 EMU == 516.1923pt). Width is
# untouched by the recompute, so it's left alone.
$shp.Height = 516.19226472

# --- 2. TLS bullet: bold + red "a successor of SSL 3.0" ------------------
# Re-read the text since it shifted slightly after the PGP edit above.
$full = $tr.Text
$highlight = "a successor of SSL 3.0"
$idxTls = $full.IndexOf($highlight)
$rngTls = $tr.Characters($idxTls + 1, $highlight.Length)
$rngTls.Font.Bold = $true
$rngTls.Font.Color.RGB = 255

# --- 3. New textbox: "SSL was deprecated in 2015" -------------------------
$inchTxt = $s.Shapes.AddTextbox(1, 694.08, 263.7853543307087, 202.83425196850393, 24.234409448818898)
$inchTxt.Fill.Visible = 0
$inchTxt.TextFrame.AutoSize = 1
$inchTxt.TextFrame.WordWrap = -1

$newTr = $inchTxt.TextFrame.TextRange
$newTr.Text = "SSL was deprecated in 2015"
$newTr.Font.Bold = $true
$newTr.Font.Color.RGB = 255

# AutoSize recalculated the height against the default live font metrics;
# pin it back to the author's recorded size (~24.23pt == 307777 EMU).
$inchTxt.Width = 202.83425196850393
$inchTxt.Height = 24.23445
